$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "randu"
$ws.Range("B2").Value = "condori"
$ws.Range("A3").Value = "randu"
$ws.Range("B3").Value = "condori"
